$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '21.304.71'
$ws.Range("E2").Value = '  +4.10%  '
$ws.Range("D3").Value = '1.546.74'
$ws.Range("E3").Value = '  +4.91%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9652'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '282.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3638'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3209'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '40.89'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.66%  '
$ws.Range("E10").Value = '  +5.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06907'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.725'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.404'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001053'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9654'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").Value = '1.543.34'
$ws.Range("E18").Value = '  +4.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06140'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.753'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.330'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.65%  '
$ws.Range("D25").Value = '21.346.36'
$ws.Range("E25").Value = '  +4.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '147.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.258'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.83'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.38%  '
$ws.Range("D29").Value = '1.712.64'
$ws.Range("E29").Value = '  +5.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.80'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.026'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.51%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8660'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.97%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.278'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08058'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.519'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.989'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.204'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05889'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02126'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.884'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1936'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9645'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5520'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.21%  '
$ws.Range("E45").Value = '  +4.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.584'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5484'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.887'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06620'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '69.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.90%  '
